$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000" or "0.08140"
# keep their exact formatting instead of being coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) values for rows 2-42 (coin order unchanged in this range) ---
$ws.Range("D2").Value = '25.779.79'
$ws.Range("E2").Value = '  -2.64%  '
$ws.Range("D3").Value = '1.745.75'
$ws.Range("E3").Value = '  -4.98%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '238.41'
$ws.Range("E5").Value = '  -9.00%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.5058'
$ws.Range("E7").Value = '  -6.17%  '
$ws.Range("D8").Value = '42.02'
$ws.Range("E8").Value = '  -6.30%  '
$ws.Range("D9").Value = '0.2730'
$ws.Range("E9").Value = '  -9.34%  '
$ws.Range("D10").Value = '0.06155'
$ws.Range("E10").Value = '  -10.96%  '
$ws.Range("D11").Value = '1.747.13'
$ws.Range("E11").Value = '  -4.87%  '
$ws.Range("D12").Value = '0.06931'
$ws.Range("E12").Value = '  -3.03%  '
$ws.Range("D13").Value = '15.52'
$ws.Range("E13").Value = '  -11.93%  '
$ws.Range("D14").Value = '4.527'
$ws.Range("E14").Value = '  -9.32%  '
$ws.Range("D15").Value = '0.6016'
$ws.Range("E15").Value = '  -18.49%  '
$ws.Range("D16").Value = '77.04'
$ws.Range("E16").Value = '  -13.69%  '
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = '25.784.59'
$ws.Range("E19").Value = '  -2.73%  '
$ws.Range("D20").Value = '0.000006894'
$ws.Range("E20").Value = '  -12.70%  '
$ws.Range("D21").Value = '11.66'
$ws.Range("E21").Value = '  -15.81%  '
$ws.Range("D22").Value = '1.968.80'
$ws.Range("E22").Value = '  -4.99%  '
$ws.Range("D23").Value = '4.065'
$ws.Range("E23").Value = '  -11.35%  '
$ws.Range("D24").Value = '5.245'
$ws.Range("E24").Value = '  -12.40%  '
$ws.Range("D25").Value = '8.176'
$ws.Range("E25").Value = '  -11.22%  '
$ws.Range("D26").Value = '137.89'
$ws.Range("E26").Value = '  -3.40%  '
$ws.Range("D27").Value = '1.467'
$ws.Range("E27").Value = '  -14.57%  '
$ws.Range("D28").Value = '1.821'
$ws.Range("E28").Value = '  -16.37%  '
$ws.Range("D29").Value = '15.01'
$ws.Range("E29").Value = '  -11.78%  '
$ws.Range("D30").Value = '103.92'
$ws.Range("E30").Value = '  -6.38%  '
$ws.Range("D31").Value = '0.08140'
$ws.Range("E31").Value = '  -8.06%  '
$ws.Range("D32").Value = '3.714'
$ws.Range("E32").Value = '  -12.61%  '
$ws.Range("D33").Value = '3.488'
$ws.Range("E33").Value = '  -13.90%  '
$ws.Range("D34").Value = '0.04544'
$ws.Range("E34").Value = '  -6.17%  '
$ws.Range("D35").Value = '0.9992'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '2.616'
$ws.Range("E36").Value = '  -10.76%  '
$ws.Range("D37").Value = '0.9857'
$ws.Range("E37").Value = '  -12.91%  '
$ws.Range("D38").Value = '0.6106'
$ws.Range("E38").Value = '  -16.43%  '
$ws.Range("D39").Value = '2.673'
$ws.Range("E39").Value = '  -13.55%  '
$ws.Range("D40").Value = '0.01556'
$ws.Range("E40").Value = '  -9.51%  '
$ws.Range("D41").Value = '1.933'
$ws.Range("E41").Value = '  -14.59%  '
$ws.Range("D42").Value = '0.9999'
$ws.Range("E42").Value = '  +0.00%  '

# --- Rows 43-51: a new "PaxosStandard" entry is inserted at row 43, shifting the existing rows down by one; ---
# --- the previous last row (Aave) falls outside the A1:E51 sheet range and is dropped. ---
$ws.Range("B43").Value = 'PaxosStandard'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '101.99'
$ws.Range("E44").Value = '  -5.65%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3851'
$ws.Range("E45").Value = '  -18.29%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '0.7407'
$ws.Range("E46").Value = '  -18.24%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '4.962'
$ws.Range("E47").Value = '  -15.93%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05383'
$ws.Range("E48").Value = '  -6.62%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1113'
$ws.Range("E49").Value = '  -11.15%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = '5.974'
$ws.Range("E50").Value = '  -19.34%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '30.18'
$ws.Range("E51").Value = '  -13.38%  '
